$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date/time value on row 99 (A99): 45484.6250347222 -> 45484.2916666667
$ws.Range("A99").Value = 45484.2916666667

# Add new row 100 with the R script results
$ws.Range("A100").Value = 45485.3740740741
$ws.Range("A100").NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Range("A100").Font.Name = "Calibri"

$ws.Range("B100").Value = 3750
$ws.Range("C100").Value = 3.25
$ws.Range("D100").Value = 3.16000008583069
$ws.Range("E100").Value = 3.16000008583069
$ws.Range("F100").Value = 3.25

# G100 must stay text ("3.25"), so force text formatting before assigning,
# then clear the formatting again so the cell keeps the default (unstyled) look
$ws.Range("G100").NumberFormat = "@"
$ws.Range("G100").Value = "3.25"
$ws.Range("G100").ClearFormats()

$ws.Range("H100").Value = "ESPE.MI"
